# Apply "User data with filename fix" changes to oda-BG workbook.
$wb = $excel.ActiveWorkbook
$wsNotes = $wb.Worksheets.Item("Notes")
$wsData = $wb.Worksheets.Item("Data")

# --- Notes sheet: fix the "Units of measure" line (row 3) ---
$wsNotes.Range("A3").Value = "Units of measure: constant 2015 US$"

# --- Data sheet: populate the previously-empty data rows (2-15) ---
$rows = @(
    @{ Id = "AF";     Name = "Afghanistan";            Year = 2015; Value = 30000 },
    @{ Id = "AL";     Name = "Albania";                Year = 2015; Value = 170000 },
    @{ Id = "BA";     Name = "Bosnia & Herzegovina";   Year = 2015; Value = 120000 },
    @{ Id = "europe"; Name = "Europe, regional";       Year = 2015; Value = 140000 },
    @{ Id = "GE";     Name = "Georgia";                Year = 2015; Value = 90000 },
    @{ Id = "IQ";     Name = "Iraq";                   Year = 2015; Value = 30000 },
    @{ Id = "MK";     Name = "Macedonia";               Year = 2015; Value = 110000 },
    @{ Id = "MD";     Name = "Moldova";                Year = 2015; Value = 310000 },
    @{ Id = "NP";     Name = "Nepal";                  Year = 2015; Value = 40000 },
    @{ Id = "RS";     Name = "Serbia";                 Year = 2015; Value = 90000 },
    @{ Id = "SO";     Name = "Somalia";                Year = 2015; Value = 30000 },
    @{ Id = "SY";     Name = "Syria";                  Year = 2015; Value = 110000 },
    @{ Id = "UA";     Name = "Ukraine";                Year = 2015; Value = 120000 },
    @{ Id = "VN";     Name = "Viet Nam";               Year = 2015; Value = 10000 }
)

$r = 2
foreach ($row in $rows) {
    $wsData.Cells.Item($r, 1).Value = $row.Id
    $wsData.Cells.Item($r, 2).Value = $row.Name
    $wsData.Cells.Item($r, 3).Value = $row.Year
    $wsData.Cells.Item($r, 4).Value = $row.Value
    $r++
}
